# Fix typos in "Bookworm Library Interface.pptx"
#  1. Date placeholder field text on the slide master + every slide layout:
#       "02/12/2014" -> "12/2/14"
#  2. Slide 10 title: "Users the have checked out..." -> "Users that have checked out..."
#  3. Slide 14 query textbox: re-join "Query " + ":" runs (no text change)
#  4. Slide 7 query textbox: re-join three split runs in three different spots
#     (no text change)

$p = $ppt.ActivePresentation

function Set-RangeText {
    param($TextRange, [int]$Start, [int]$Length, [string]$NewText)
    $sub = $TextRange.Characters($Start, $Length)
    $sub.Text = $NewText
}

function Update-DatePlaceholder {
    param($Shapes, [string]$NewText)
    for ($j = 1; $j -le $Shapes.Count; $j++) {
        $shp = $Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            $len = $tr.Text.Length
            if ($len -gt 0) {
                Set-RangeText $tr 1 $len $NewText
            }
            break
        }
    }
}

# --- 1. slide master + every slide layout date field -----------------------
Update-DatePlaceholder $p.SlideMaster.Shapes "12/2/14"

$layouts = $p.Designs.Item(1).SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i).Shapes "12/2/14"
}

# --- 2. slide 10 title typo: "the" -> "that" --------------------------------
$s10 = $p.Slides.Item(10)
$title10 = $s10.Shapes.Item(1).TextFrame.TextRange
Set-RangeText $title10 7 4 "that "

# --- 3. slide 14: merge "Query " + ":" into a single run --------------------
$s14 = $p.Slides.Item(14)
$tb14 = $s14.Shapes.Item(2).TextFrame.TextRange
Set-RangeText $tb14 1 7 "Query :"

# --- 4. slide 7: merge three split runs back together -----------------------
$s7 = $p.Slides.Item(7)
$tb7 = $s7.Shapes.Item(2).TextFrame.TextRange

Set-RangeText $tb7 58 39 ", `tr1.publisher_id, Publishers.name as "
Set-RangeText $tb7 158 39 "`tFROM Books) as r1, Authors, Publishers"
Set-RangeText $tb7 234 22 " AND `tr1.publisher_id="
